# Updated cryptos list on Wed Oct  2 05:07:45 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds numeric-looking text (e.g. "554.15", "61.542.75").
# Excel's COM layer auto-converts plain numeric-looking strings assigned via
# .Value into real numbers, which would lose the original text formatting
# (e.g. thousands separated by dots like "61.542.75") and introduce float
# rounding noise. Forcing the cell to Text format ("@") before assignment
# keeps these values as literal strings, matching the source data.
$priceCells = @("D2","D3","D5","D6","D8","D9","D11","D12","D14","D15","D17","D18","D19","D20","D22","D24","D25","D26","D27","D30","D31","D32","D36","D40","D41","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "61.542.75"
$ws.Range("E2").Value = "  -3.57%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.488.77"
$ws.Range("E3").Value = "  -5.79%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.03%  "

# Row 5 - BNB
$ws.Range("D5").Value = "554.15"
$ws.Range("E5").Value = "  -4.60%  "

# Row 6 - Solana
$ws.Range("D6").Value = "146.80"
$ws.Range("E6").Value = "  -5.68%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.03%  "

# Row 8 - XRP
$ws.Range("D8").Value = "0.604"
$ws.Range("E8").Value = "  -2.48%  "

# Row 9 - LidoStakedEther
$ws.Range("D9").Value = "2.486.84"
$ws.Range("E9").Value = "  -5.79%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -8.75%  "

# Row 11 - was TRON, now Toncoin
$ws.Range("B11").Value = "Toncoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D11").Value = "5.44"
$ws.Range("E11").Value = "  -6.45%  "

# Row 12 - was Toncoin, now TRON
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "0.154"
$ws.Range("E12").Value = "  -1.49%  "

# Row 13 - Cardano
$ws.Range("E13").Value = "  -6.37%  "

# Row 14 - Avalanche
$ws.Range("D14").Value = "26.29"
$ws.Range("E14").Value = "  -7.71%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "2.935.26"
$ws.Range("E15").Value = "  -5.94%  "

# Row 16 - ShibaInu
$ws.Range("E16").Value = "  -8.71%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "61.484.25"
$ws.Range("E17").Value = "  -3.67%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "2.486.84"
$ws.Range("E18").Value = "  -5.61%  "

# Row 19 - Chainlink
$ws.Range("D19").Value = "11.19"
$ws.Range("E19").Value = "  -7.90%  "

# Row 20 - Uniswap
$ws.Range("D20").Value = "7.01"
$ws.Range("E20").Value = "  -8.48%  "

# Row 21 - Polkadot
$ws.Range("E21").Value = "  -7.07%  "

# Row 22 - BitcoinCash
$ws.Range("D22").Value = "322.71"
$ws.Range("E22").Value = "  -6.39%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  -0.03%  "

# Row 24 - SuiNetwork
$ws.Range("D24").Value = "1.84"
$ws.Range("E24").Value = "  -2.23%  "

# Row 25 - Litecoin
$ws.Range("D25").Value = "64.04"
$ws.Range("E25").Value = "  -5.99%  "

# Row 26 - PEPE
$ws.Range("D26").Value = "0.0₃0994"
$ws.Range("E26").Value = "  -9.05%  "

# Row 27 - WrappedeETH
$ws.Range("D27").Value = "2.607.00"
$ws.Range("E27").Value = "  -5.79%  "

# Row 28 - Fetch.AI
$ws.Range("E28").Value = "  -5.63%  "

# Row 29 - Binance-PegBSC-USD
$ws.Range("E29").Value = "  +0.09%  "

# Row 30 - Bittensor
$ws.Range("D30").Value = "539.81"
$ws.Range("E30").Value = "  -10.69%  "

# Row 31 - InternetComputer(DFINITY)
$ws.Range("D31").Value = "8.36"
$ws.Range("E31").Value = "  -9.84%  "

# Row 32 - Aptos
$ws.Range("D32").Value = "7.70"
$ws.Range("E32").Value = "  -5.20%  "

# Row 33 - Kaspa
$ws.Range("E33").Value = "  -6.06%  "

# Row 34 - PancakeSwap
$ws.Range("E34").Value = "  -7.34%  "

# Row 35 - ImmutableX
$ws.Range("E35").Value = "  -8.82%  "

# Row 36 - RenderToken
$ws.Range("D36").Value = "5.91"
$ws.Range("E36").Value = "  -10.39%  "

# Row 37 - NEARProtocol
$ws.Range("E37").Value = "  -10.41%  "

# Row 38 - FirstDigitalUSD
$ws.Range("E38").Value = "  -0.12%  "

# Row 39 - PolygonEcosystemToken
$ws.Range("E39").Value = "  -5.50%  "

# Row 40 - EthereumClassic
$ws.Range("D40").Value = "18.57"
$ws.Range("E40").Value = "  -6.04%  "

# Row 41 - Monero
$ws.Range("D41").Value = "148.84"
$ws.Range("E41").Value = "  -1.59%  "

# Row 42 - Stacks
$ws.Range("E42").Value = "  -8.66%  "

# Row 43 - USDe
$ws.Range("E43").Value = "  +0.09%  "

# Row 44 - OKB
$ws.Range("D44").Value = "40.40"
$ws.Range("E44").Value = "  -3.59%  "

# Row 45 - dogwifhat
$ws.Range("D45").Value = "2.37"
$ws.Range("E45").Value = "  -7.50%  "

# Row 46 - Aave
$ws.Range("D46").Value = "148.34"
$ws.Range("E46").Value = "  -7.40%  "

# Row 47 - Filecoin
$ws.Range("D47").Value = "3.64"
$ws.Range("E47").Value = "  -6.86%  "

# Row 48 - InjectiveProtocol
$ws.Range("D48").Value = "21.25"
$ws.Range("E48").Value = "  -12.89%  "

# Row 49 - Hedera
$ws.Range("D49").Value = "0.0539"
$ws.Range("E49").Value = "  -8.20%  "

# Row 50 - Mantle
$ws.Range("D50").Value = "0.598"
$ws.Range("E50").Value = "  -5.66%  "

# Row 51 - Stellar
$ws.Range("D51").Value = "0.0948"
$ws.Range("E51").Value = "  -4.87%  "
